# Add more Aviation Ground-staff KPIs to the "Aerospace" sheet (Domains and KPI's.xlsx)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Aerospace")
$ws.Activate()

# --- New data rows (13-19) -------------------------------------------------
# Row 13: new Category "Ground Staff Maintenance" + first KPI
$ws.Range("A13").Value = "Ground Staff Maintenance"
$ws.Range("B13").Value = "Turnaround Time per Aircraft Type"
$ws.Range("C13").Value = "Total time elapsed between block-on and block-off. This metric must consider aircraft configuration, gate location, and local constraints."

$ws.Range("B14").Value = "Average Delay per Flight"
$ws.Range("C14").Value = "Captures the mean delay duration across all flights. Used in conjunction with root-cause analysis, it helps isolate the share of delay attributable to ground handling."

$ws.Range("B15").Value = "Mishandled Baggage Index (MBI)"
$ws.Range("C15").Value = "Measures the percentage of bags delivered to the correct flight or carousel within SLA thresholds."

$ws.Range("B16").Value = "Delays per 100 Flight Legs"
$ws.Range("C16").Value = "Quantifies how frequently delays occur across a standardized number of flights, helping detect systemic issues in scheduling or execution."

$ws.Range("B17").Value = "Ground Incident Rate per 1,000 Flights"
$ws.Range("C17").Value = "Tracks safety-related deviations, including equipment collisions, misconnects, and FOD events."

$ws.Range("B18").Value = "Staff Task Compliance"
$ws.Range("C18").Value = "Monitors procedural adherence based on logged task completions, often integrated with mobile crew management systems."

$ws.Range("B19").Value = "Passenger Touchpoint SLA Compliance"
$ws.Range("C19").Value = "(e.g., wheelchair delivery, cabin cleaning readiness) connects service-level commitments to operational execution."

# --- Formatting: wrap the "Short Description" column like the rest of the sheet
$ws.Range("C13:C19").WrapText = $true

# --- Row heights (matches the two-line wrapped rows elsewhere on the sheet)
$ws.Rows.Item(13).RowHeight = 30
$ws.Rows.Item(14).RowHeight = 30
$ws.Rows.Item(16).RowHeight = 30
$ws.Rows.Item(18).RowHeight = 30

# --- Column widths: column A is now used, column B got noticeably wider
$ws.Columns.Item(1).ColumnWidth = 23.6
$ws.Columns.Item(2).ColumnWidth = 36

# --- View state: scroll so the new rows are visible, select C22 like the source file
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C22").Select()
